# This script applies a weekly reshuffle of the Alcachofa price records
# (rows 2-17) in the active worksheet. Each data row's Fecha/Variedad/
# Calidad/Volumen/Precio/Unidad fields are updated in place so the overall
# set of records matches the latest weekly export; row 14 is left untouched
# because its figures did not change in the new export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2 (now holds the data previously in row 12)
$ws.Range("D2").Value = 44400
$ws.Range("H2").Value = 'Española'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("N2").Value = '$/caja 30 unidades'
$ws.Range("P2").Value = 500
$ws.Range("Q2").Value = 30

# Row 3 (now holds the data previously in row 17)
$ws.Range("D3").Value = 44421
$ws.Range("H3").Value = 'Española'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 16500
$ws.Range("L3").Value = 16500
$ws.Range("M3").Value = 16500
$ws.Range("N3").Value = '$/caja 30 unidades'
$ws.Range("P3").Value = 550
$ws.Range("Q3").Value = 30

# Row 4 (now holds the data previously in row 10)
$ws.Range("D4").Value = 44176
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("N4").Value = '$/caja 40 unidades'
$ws.Range("P4").Value = 275
$ws.Range("Q4").Value = 40

# Row 5 (now holds the data previously in row 7)
$ws.Range("D5").Value = 44484
$ws.Range("H5").Value = 'Madrigal'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 110
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("N5").Value = '$/caja 50 unidades'
$ws.Range("P5").Value = 220
$ws.Range("Q5").Value = 50

# Row 6 (now holds the data previously in row 2)
$ws.Range("D6").Value = 44488
$ws.Range("H6").Value = 'Madrigal'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 12000
$ws.Range("N6").Value = '$/caja 40 unidades'
$ws.Range("P6").Value = 300
$ws.Range("Q6").Value = 40

# Row 7 (now holds the data previously in row 8)
$ws.Range("D7").Value = 44386
$ws.Range("H7").Value = 'Española'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("P7").Value = 500
$ws.Range("Q7").Value = 30

# Row 8 (now holds the data previously in row 11)
$ws.Range("D8").Value = 44390
$ws.Range("H8").Value = 'Española'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 16000
$ws.Range("N8").Value = '$/caja 30 unidades'
$ws.Range("P8").Value = 533
$ws.Range("Q8").Value = 30

# Row 9 (now holds the data previously in row 16)
$ws.Range("D9").Value = 44166
$ws.Range("H9").Value = 'Madrigal'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("P9").Value = 250
$ws.Range("Q9").Value = 40

# Row 10 (now holds the data previously in row 5)
$ws.Range("D10").Value = 44481
$ws.Range("H10").Value = 'Madrigal'
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 11000
$ws.Range("N10").Value = '$/caja 50 unidades'
$ws.Range("P10").Value = 220
$ws.Range("Q10").Value = 50

# Row 11 (now holds the data previously in row 3)
$ws.Range("D11").Value = 44491
$ws.Range("H11").Value = 'Madrigal'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 11000
$ws.Range("M11").Value = 11000
$ws.Range("N11").Value = '$/caja 40 unidades'
$ws.Range("P11").Value = 275
$ws.Range("Q11").Value = 40

# Row 12 (now holds the data previously in row 4)
$ws.Range("D12").Value = 44162
$ws.Range("H12").Value = 'Madrigal'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = '$/caja 40 unidades'
$ws.Range("P12").Value = 250
$ws.Range("Q12").Value = 40

# Row 13 (now holds the data previously in row 9)
$ws.Range("D13").Value = 44446
$ws.Range("H13").Value = 'Madrigal'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 120
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 16000
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("P13").Value = 400
$ws.Range("Q13").Value = 40

# Row 14: unchanged

# Row 15 (now holds the data previously in row 6)
$ws.Range("D15").Value = 44161
$ws.Range("H15").Value = 'Madrigal'
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 11000
$ws.Range("M15").Value = 11000
$ws.Range("N15").Value = '$/caja 40 unidades'
$ws.Range("P15").Value = 275
$ws.Range("Q15").Value = 40

# Row 16 (now holds the data previously in row 15)
$ws.Range("D16").Value = 44407
$ws.Range("H16").Value = 'Española'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("N16").Value = '$/caja 30 unidades'
$ws.Range("P16").Value = 600
$ws.Range("Q16").Value = 30

# Row 17 (now holds the data previously in row 13)
$ws.Range("D17").Value = 44418
$ws.Range("H17").Value = 'Española'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 16000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 16000
$ws.Range("N17").Value = '$/caja 30 unidades'
$ws.Range("P17").Value = 533
$ws.Range("Q17").Value = 30
